$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.035.81'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '2.370.53'
$ws.Range("E3").Value = '  +1.76%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.502'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.482'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.32%  '
$ws.Range("E11").Value = '  +2.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0786'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.75'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.05%  '
$ws.Range("D15").Value = '2.730.35'
$ws.Range("E15").Value = '  +1.31%  '
$ws.Range("D16").Value = '2.346.02'
$ws.Range("E16").Value = '  +3.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.798'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").Value = '42.994.51'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.64%  '
$ws.Range("D21").Value = '0.0₃0886'
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.78%  '
$ws.Range("E24").Value = '  -2.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.73%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.49%  '
$ws.Range("E28").Value = '  +15.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.58'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0716'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.88%  '
$ws.Range("E35").Value = '  +4.13%  '
$ws.Range("E36").Value = '  +2.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.70%  '
$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.26'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.05%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '122.29'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -10.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.108'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.91%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.94%  '
$ws.Range("D43").Value = '1.934.59'
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0279'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("E45").Value = '  +4.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.31'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("D48").Value = '2.592.69'
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.77'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.13'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.72%  '
